$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Replace the two summary columns ("Completed Courses" / "Uncompleted
# Courses") with one column per individual course (F1:M1). Copy the
# formatting from the existing "Category" header (E1) first so the new
# header cells pick up the same bold / bordered / centered style instead
# of a brand new one being created.
$ws.Range("E1").Copy()
$ws.Range("F1:M1").PasteSpecial(-4122)

$ws.Range("F1").Value = "DHA Accommodations (1 hr)"
$ws.Range("G1").Value = "Leadership Training (4 hrs)"
$ws.Range("H1").Value = "MHS Customer Service (1 hr)"
$ws.Range("I1").Value = "Counterintelligence (1 hr)"
$ws.Range("J1").Value = "HIPAA Training (1 hr)"
$ws.Range("K1").Value = "Supervisor Safety Training (2 hrs)"
$ws.Range("L1").Value = "Employee Safety (1 hr)"
$ws.Range("M1").Value = "Violence Response (1 hr)"

# Helper: leave a cell blank (no completion status for that person/course)
# while still keeping the cell present with the plain/default style so the
# sheet keeps a full rectangular A:M grid on every data row.
function Clear-Status($rng) {
    $rng.ClearContents()
    $rng.Borders.LineStyle = -4142
}

# --- Data rows ----------------------------------------------------------
# The old F/G columns held Python-list-literal strings of completed /
# uncompleted course names. They're replaced by a Completed / NOT
# Completed / blank flag per course column (F:M).

# Row 2 - John Doe
$ws.Range("F2").Value = "Completed"
$ws.Range("G2").Value = "Completed"
$ws.Range("H2").Value = "Completed"
Clear-Status $ws.Range("I2")
Clear-Status $ws.Range("J2")
Clear-Status $ws.Range("K2")
Clear-Status $ws.Range("L2")
Clear-Status $ws.Range("M2")

# Row 3 - Andrew Hartmann
Clear-Status $ws.Range("F3")
Clear-Status $ws.Range("G3")
Clear-Status $ws.Range("H3")
$ws.Range("I3").Value = "NOT Completed"
$ws.Range("J3").Value = "NOT Completed"
$ws.Range("K3").Value = "Completed"
$ws.Range("L3").Value = "Completed"
Clear-Status $ws.Range("M3")

# Row 4 - Nick Fletcher
Clear-Status $ws.Range("F4")
Clear-Status $ws.Range("G4")
Clear-Status $ws.Range("H4")
$ws.Range("I4").Value = "Completed"
$ws.Range("J4").Value = "Completed"
$ws.Range("K4").Value = "Completed"
$ws.Range("L4").Value = "Completed"
Clear-Status $ws.Range("M4")

# Row 5 - John Cena
$ws.Range("F5").Value = "Completed"
$ws.Range("G5").Value = "Completed"
$ws.Range("H5").Value = "Completed"
Clear-Status $ws.Range("I5")
Clear-Status $ws.Range("J5")
Clear-Status $ws.Range("K5")
Clear-Status $ws.Range("L5")
$ws.Range("M5").Value = "Completed"
